# Updated symbol list on Tue Jan 31 08:49:10 UTC 2023 with GitHub Actions
#
# Refreshes the crypto price/volume snapshot (and, for two rows, the coin
# identity itself) on the active worksheet to match the latest pull from
# coinranking.com. Price/volume columns are stored as plain text in this
# sheet (so values like "0.0002000" keep their exact trailing zeros and
# "-0.21%" stays literal text instead of becoming a formatted number), so
# each target cell is pre-formatted as Text before the new value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '310.72' },
    @{ Cell = 'E2'; Value = '-0.21%' },
    @{ Cell = 'D3'; Value = '37.58' },
    @{ Cell = 'E3'; Value = '-1.56%' },
    @{ Cell = 'D4'; Value = '5.072' },
    @{ Cell = 'E4'; Value = '-1.27%' },
    @{ Cell = 'D5'; Value = '0.07763' },
    @{ Cell = 'E5'; Value = '-4.14%' },
    @{ Cell = 'D6'; Value = '4.353' },
    @{ Cell = 'E6'; Value = '-2.13%' },
    @{ Cell = 'D7'; Value = '8.224' },
    @{ Cell = 'E7'; Value = '-0.92%' },
    @{ Cell = 'D8'; Value = '1.881' },
    @{ Cell = 'E8'; Value = '-3.49%' },
    @{ Cell = 'D9'; Value = '0.9212' },
    @{ Cell = 'E10'; Value = '-8.61%' },
    @{ Cell = 'D11'; Value = '0.1923' },
    @{ Cell = 'E11'; Value = '-1.62%' },
    @{ Cell = 'D12'; Value = '0.09334' },
    @{ Cell = 'E12'; Value = '3.19%' },
    @{ Cell = 'D13'; Value = '0.03434' },
    @{ Cell = 'E13'; Value = '-1.64%' },
    @{ Cell = 'D14'; Value = '0.09674' },
    @{ Cell = 'E14'; Value = '-0.18%' },
    @{ Cell = 'D15'; Value = '0.001371' },
    @{ Cell = 'E15'; Value = '-2.63%' },
    @{ Cell = 'D16'; Value = '0.005824' },
    @{ Cell = 'E16'; Value = '-1.05%' },
    @{ Cell = 'D17'; Value = '3.554' },
    @{ Cell = 'E17'; Value = '0.07%' },
    @{ Cell = 'E18'; Value = '-10.78%' },
    @{ Cell = 'D19'; Value = '0.3401' },
    @{ Cell = 'E19'; Value = '-1.87%' },
    @{ Cell = 'D20'; Value = '5.306' },
    @{ Cell = 'E20'; Value = '5.71%' },
    @{ Cell = 'D21'; Value = '0.1297' },
    @{ Cell = 'E21'; Value = '1.05%' },
    @{ Cell = 'E23'; Value = '5,588.15%' },
    @{ Cell = 'D24'; Value = '0.04357' },
    @{ Cell = 'E24'; Value = '-0.30%' },
    @{ Cell = 'D25'; Value = '0.001212' },
    @{ Cell = 'E25'; Value = '-2.17%' },
    @{ Cell = 'D26'; Value = '0.004254' },
    @{ Cell = 'E26'; Value = '-9.98%' },
    @{ Cell = 'D27'; Value = '0.0001301' },
    @{ Cell = 'E27'; Value = '-66.23%' },
    @{ Cell = 'D39'; Value = '0.02091' },
    @{ Cell = 'E39'; Value = '-5.27%' },
    @{ Cell = 'D40'; Value = '0.04978' },
    @{ Cell = 'E40'; Value = '-4.90%' },
    @{ Cell = 'D41'; Value = '0.007640' },
    @{ Cell = 'E41'; Value = '0.68%' },
    @{ Cell = 'D42'; Value = '0.009830' },
    @{ Cell = 'E42'; Value = '-4.70%' },
    @{ Cell = 'D43'; Value = '0.1343' },
    @{ Cell = 'E43'; Value = '-3.41%' },
    @{ Cell = 'E44'; Value = '-2.16%' },
    @{ Cell = 'D45'; Value = '0.008818' },
    @{ Cell = 'E45'; Value = '-3.21%' },
    @{ Cell = 'D46'; Value = '0.00006649' },
    @{ Cell = 'E46'; Value = '0.44%' },
    @{ Cell = 'E47'; Value = '-0.24%' },
    @{ Cell = 'B48'; Value = 'BOLO' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo' },
    @{ Cell = 'D48'; Value = '0.002936' },
    @{ Cell = 'E48'; Value = '-2.60%' },
    @{ Cell = 'B49'; Value = 'CoinbaseStockToken' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin' },
    @{ Cell = 'D49'; Value = '0.001201' },
    @{ Cell = 'E49'; Value = '-29.01%' },
    @{ Cell = 'E50'; Value = '-0.24%' },
    @{ Cell = 'D51'; Value = '0.0002001' },
    @{ Cell = 'E51'; Value = '-0.24%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Columns D (Price) and E (Volume(1h)) hold numeric-looking text; force
    # Text formatting first so Excel doesn't silently coerce "310.72" or
    # "-0.21%" into a Double. Columns B/C (Coin/Link) are plain text already.
    if ($u.Cell -match '^[DE]') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates"
